$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("individual")
$ws.Range("I2").Value = "1.0"
